# Germanize the Assignment4 workbook: rename sheets, translate headers,
# fix the March sheet header cell styles, and update the active
# tab / selections to match the reviewed state.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the month sheets (English -> German) ---------------------
$wsJan = $wb.Worksheets.Item(1)   # "January"
$wsFeb = $wb.Worksheets.Item(2)   # "February"
$wsMar = $wb.Worksheets.Item(4)   # "March"

$wsJan.Name = "Januar"
$wsFeb.Name = "Februar"
$wsMar.Name = "März"

# --- 2. Translate the header rows on each roster sheet -------------------
# Row 1: title "Webinar Anwesenheitsliste: <Month>"
# Row 2: column headers Email Addresse / Vorname / Nachname / Webinar Bestanden:

$wsJan.Range("A1").Value = "Webinar Anwesenheitsliste: Januar"
$wsJan.Range("A2").Value = "Email Addresse"
$wsJan.Range("B2").Value = "Vorname"
$wsJan.Range("C2").Value = "Nachname"
$wsJan.Range("D2").Value = "Webinar Bestanden:"

$wsFeb.Range("A1").Value = "Webinar Anwesenheitsliste: Februar"
$wsFeb.Range("A2").Value = "Email Addresse"
$wsFeb.Range("B2").Value = "Vorname"
$wsFeb.Range("C2").Value = "Nachname"
$wsFeb.Range("D2").Value = "Webinar Bestanden:"

$wsMar.Range("A1").Value = "Webinar Anwesenheitsliste: März"
$wsMar.Range("A2").Value = "Email Addresse"
$wsMar.Range("B2").Value = "Vorname"
$wsMar.Range("C2").Value = "Nachname"
$wsMar.Range("D2").Value = "Webinar Bestanden:"

# The March sheet's D1/D2 header cells carried a stray date number format
# left over from an older layout; bring them back in line with the rest
# of the header row (copy A1/A2's format onto them).
$wsMar.Range("A1").Copy()
$wsMar.Range("D1").PasteSpecial(-4122)
$wsMar.Range("A2").Copy()
$wsMar.Range("D2").PasteSpecial(-4122)

# --- 3. Update tab order / selection state --------------------------------
# NOTE: Range.Select() implicitly activates the sheet it's called on, so
# the calls below are ordered such that Februar (the desired active tab)
# is both activated and selected *last*.

# März keeps its own selection independent of which tab ends up active.
$wsMar.Range("A2").Select()

# Januar is no longer the active tab; just leave a range selected on it.
$wsJan.Activate()
$wsJan.Range("A1:D2").Select()

# Februar becomes the active tab.
$wsFeb.Activate()
$wsFeb.Range("A7").Select()
